$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9261095
$ws.Range("J17").Value = 9261095
$ws.Range("L17").Value = 27783285
$ws.Range("N17").Value = -27783621
$ws.Range("H33").Value = 171.5
$ws.Range("I33").Value = 174.8
$ws.Range("J33").Value = 155
$ws.Range("K33").Value = 174.8
$ws.Range("L33").Value = 155
$ws.Range("M33").Value = 54.19999999999999
$ws.Range("N33").Value = -613
$ws.Range("H132").Value = 5705.2666
$ws.Range("I132").Value = 5897.6895
$ws.Range("J132").Value = 125
$ws.Range("K132").Value = 17693.0685
$ws.Range("L132").Value = 375
$ws.Range("M132").Value = -15163.0685
$ws.Range("N132").Value = -5435
$ws.Range("H137").Value = 5203.8887
$ws.Range("I137").Value = 1662.9375
$ws.Range("J137").Value = 10354.363
$ws.Range("K137").Value = 4988.8125
$ws.Range("L137").Value = 31063.089
$ws.Range("M137").Value = -2438.8125
$ws.Range("N137").Value = -36163.089
$ws.Range("H138").Value = 304895.7
$ws.Range("I138").Value = 4242.125
$ws.Range("J138").Value = 427195.44
$ws.Range("K138").Value = 12726.375
$ws.Range("L138").Value = 1281586.32
$ws.Range("M138").Value = -7586.375
$ws.Range("N138").Value = -1291866.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3250
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H32").Value = 3653.3699
$ws.Range("I32").Value = 2894.1594
$ws.Range("K32").Value = 2894.1594
$ws.Range("M32").Value = -2607.1594
$ws.Range("H45").Value = 30516.445
$ws.Range("I45").Value = 41695.168
$ws.Range("K45").Value = 41695.168
$ws.Range("M45").Value = -41318.168
$ws.Range("H61").Value = 5037.357
$ws.Range("I61").Value = 2222.2
$ws.Range("J61").Value = 6601.3335
$ws.Range("K61").Value = 2222.2
$ws.Range("L61").Value = 6601.3335
$ws.Range("M61").Value = -2010.2
$ws.Range("N61").Value = -7025.3335
$ws.Range("H109").Value = 94999
$ws.Range("J109").Value = 94999
$ws.Range("L109").Value = 94999
$ws.Range("N109").Value = -97773
$ws.Range("H122").Value = 3051.0393
$ws.Range("I122").Value = 2662.5217
$ws.Range("K122").Value = 7987.5651
$ws.Range("M122").Value = -5537.5651
$ws.Range("H132").Value = 1929.9807
$ws.Range("I132").Value = 1325.5264
$ws.Range("K132").Value = 3976.5792
$ws.Range("M132").Value = -1446.5792
$ws.Range("H136").Value = 5037.357
$ws.Range("I136").Value = 2222.2
$ws.Range("J136").Value = 6601.3335
$ws.Range("K136").Value = 6666.599999999999
$ws.Range("L136").Value = 19804.0005
$ws.Range("M136").Value = -4116.599999999999
$ws.Range("N136").Value = -24904.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1482.8572
$ws.Range("I22").Value = 1230
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1230
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1057
$ws.Range("N22").Value = -3346
$ws.Range("H97").Value = 10467
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982
$ws.Range("H107").Value = 4050255.2
$ws.Range("I107").Value = 4809207
$ws.Range("J107").Value = 2511
$ws.Range("K107").Value = 4809207
$ws.Range("L107").Value = 2511
$ws.Range("M107").Value = -4807287
$ws.Range("N107").Value = -6351
$ws.Range("H132").Value = 84500
$ws.Range("J132").Value = 84500
$ws.Range("L132").Value = 84500
$ws.Range("N132").Value = -94620

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1817
$ws.Range("I22").Value = 1498.5
$ws.Range("J22").Value = 2199.2
$ws.Range("K22").Value = 1498.5
$ws.Range("L22").Value = 2199.2
$ws.Range("M22").Value = -1148.5
$ws.Range("N22").Value = -2899.2
$ws.Range("H58").Value = 3631.4092
$ws.Range("I58").Value = 2693.9092
$ws.Range("K58").Value = 2693.9092
$ws.Range("M58").Value = -2490.9092
$ws.Range("H62").Value = 16687666
$ws.Range("I62").Value = 25009498
$ws.Range("K62").Value = 25009498
$ws.Range("M62").Value = -25008874
$ws.Range("H65").Value = 16687666
$ws.Range("I65").Value = 25009498
$ws.Range("K65").Value = 125047490
$ws.Range("M65").Value = -125044370
$ws.Range("H86").Value = 1294011.1
$ws.Range("I86").Value = 3381.5454
$ws.Range("J86").Value = 3322143.5
$ws.Range("K86").Value = 3381.5454
$ws.Range("L86").Value = 3322143.5
$ws.Range("M86").Value = -2258.5454
$ws.Range("N86").Value = -3324389.5
$ws.Range("H89").Value = 1294011.1
$ws.Range("I89").Value = 3381.5454
$ws.Range("J89").Value = 3322143.5
$ws.Range("K89").Value = 16907.727
$ws.Range("L89").Value = 16610717.5
$ws.Range("M89").Value = -11291.727
$ws.Range("N89").Value = -16621949.5
$ws.Range("H134").Value = 5730.857
$ws.Range("I134").Value = 5519.3335
$ws.Range("K134").Value = 16558.0005
$ws.Range("M134").Value = -14023.0005
$ws.Range("H136").Value = 3631.4092
$ws.Range("I136").Value = 2693.9092
$ws.Range("K136").Value = 8081.7276
$ws.Range("M136").Value = -5531.7276

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1957.0834
$ws.Range("J2").Value = 2901.125
$ws.Range("L2").Value = 17406.75
$ws.Range("N2").Value = -17632.75
$ws.Range("H56").Value = 6905
$ws.Range("I56").Value = 6905
$ws.Range("K56").Value = 6905
$ws.Range("M56").Value = -6375
$ws.Range("H69").Value = 2006
$ws.Range("I69").Value = 512
$ws.Range("K69").Value = 1536
$ws.Range("M69").Value = -725
$ws.Range("H72").Value = 2006
$ws.Range("I72").Value = 512
$ws.Range("K72").Value = 4608
$ws.Range("M72").Value = -552
$ws.Range("H121").Value = 222527.22
$ws.Range("I121").Value = 112.666664
$ws.Range("K121").Value = 337.999992
$ws.Range("M121").Value = 972.000008
$ws.Range("H132").Value = 3400.0667
$ws.Range("I132").Value = 2249.25
$ws.Range("J132").Value = 3818.5454
$ws.Range("K132").Value = 20243.25
$ws.Range("L132").Value = 34366.9086
$ws.Range("M132").Value = -17713.25
$ws.Range("N132").Value = -39426.9086
$ws.Range("H133").Value = 7975
$ws.Range("I133").Value = 7975
$ws.Range("K133").Value = 23925
$ws.Range("M133").Value = -18865
$ws.Range("H141").Value = 31777.6
$ws.Range("I141").Value = 19444
$ws.Range("K141").Value = 58332
$ws.Range("M141").Value = -53152

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 5399.8335
$ws.Range("I31").Value = 1350
$ws.Range("J31").Value = 13499.5
$ws.Range("K31").Value = 1350
$ws.Range("L31").Value = 13499.5
$ws.Range("M31").Value = -1058
$ws.Range("N31").Value = -14083.5
$ws.Range("H37").Value = 5399.8335
$ws.Range("I37").Value = 1350
$ws.Range("J37").Value = 13499.5
$ws.Range("K37").Value = 1350
$ws.Range("L37").Value = 13499.5
$ws.Range("M37").Value = -1073
$ws.Range("N37").Value = -14053.5
$ws.Range("H70").Value = 81605.46000000001
$ws.Range("I70").Value = 129515.5
$ws.Range("K70").Value = 129515.5
$ws.Range("M70").Value = -129245.5
$ws.Range("H73").Value = 81605.46000000001
$ws.Range("I73").Value = 129515.5
$ws.Range("K73").Value = 129515.5
$ws.Range("M73").Value = -128579.5
$ws.Range("H97").Value = 2122.074
$ws.Range("I97").Value = 661.2381
$ws.Range("K97").Value = 661.2381
$ws.Range("M97").Value = -165.2381
$ws.Range("H102").Value = 4298.983
$ws.Range("I102").Value = 683.0208
$ws.Range("J102").Value = 21655.6
$ws.Range("K102").Value = 683.0208
$ws.Range("L102").Value = 21655.6
$ws.Range("M102").Value = 938.9792
$ws.Range("N102").Value = -24899.6
$ws.Range("H107").Value = 1012.4286
$ws.Range("I107").Value = 851
$ws.Range("J107").Value = 1227.6666
$ws.Range("K107").Value = 851
$ws.Range("L107").Value = 1227.6666
$ws.Range("M107").Value = 1069
$ws.Range("N107").Value = -5067.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4422.5356
$ws.Range("J7").Value = 7369.4
$ws.Range("L7").Value = 7369.4
$ws.Range("N7").Value = -7593.4
$ws.Range("H22").Value = 408
$ws.Range("I22").Value = 506.83334
$ws.Range("J22").Value = 259.75
$ws.Range("K22").Value = 506.83334
$ws.Range("L22").Value = 259.75
$ws.Range("M22").Value = -211.83334
$ws.Range("N22").Value = -849.75
$ws.Range("H27").Value = 408
$ws.Range("I27").Value = 506.83334
$ws.Range("J27").Value = 259.75
$ws.Range("K27").Value = 506.83334
$ws.Range("L27").Value = 259.75
$ws.Range("M27").Value = -399.83334
$ws.Range("N27").Value = -473.75
$ws.Range("H108").Value = 87500
$ws.Range("J108").Value = 87500
$ws.Range("L108").Value = 87500
$ws.Range("N108").Value = -95180
$ws.Range("H126").Value = 4422.5356
$ws.Range("J126").Value = 7369.4
$ws.Range("L126").Value = 22108.2
$ws.Range("N126").Value = -27048.2
$ws.Range("H132").Value = 7280.76
$ws.Range("I132").Value = 5845.4736
$ws.Range("K132").Value = 17536.4208
$ws.Range("M132").Value = -15006.4208
$ws.Range("H136").Value = 5308.65
$ws.Range("I136").Value = 5399.154
$ws.Range("K136").Value = 16197.462
$ws.Range("M136").Value = -13647.462

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13997
$ws.Range("J62").Value = 13997
$ws.Range("L62").Value = 13997
$ws.Range("N62").Value = -15245
$ws.Range("H65").Value = 13997
$ws.Range("J65").Value = 13997
$ws.Range("L65").Value = 69985
$ws.Range("N65").Value = -76225
$ws.Range("H126").Value = 2199.3635
$ws.Range("I126").Value = 2020.3158
$ws.Range("K126").Value = 6060.9474
$ws.Range("M126").Value = -3590.9474
